$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'314.89"
$ws.Range("E2").Value = "'3.24%"
$ws.Range("D3").Value = "'39.39"
$ws.Range("E3").Value = "'2.81%"
$ws.Range("D4").Value = "'5.116"
$ws.Range("E4").Value = "'0.36%"
$ws.Range("D5").Value = "'0.08196"
$ws.Range("D6").Value = "'1.969"
$ws.Range("E6").Value = "'2.38%"
$ws.Range("D7").Value = "'8.243"
$ws.Range("E7").Value = "'3.81%"
$ws.Range("D8").Value = "'0.9310"
$ws.Range("E8").Value = "'0.37%"
$ws.Range("D9").Value = "'0.1405"
$ws.Range("E9").Value = "'-2.94%"
$ws.Range("D10").Value = "'0.1973"
$ws.Range("E10").Value = "'2.74%"
$ws.Range("D11").Value = "'0.09098"
$ws.Range("E11").Value = "'1.77%"
$ws.Range("D12").Value = "'0.03533"
$ws.Range("E12").Value = "'0.05%"
$ws.Range("D13").Value = "'0.09810"
$ws.Range("E13").Value = "'0.39%"
$ws.Range("D14").Value = "'0.001398"
$ws.Range("E14").Value = "'0.07%"
$ws.Range("D15").Value = "'0.005942"
$ws.Range("E15").Value = "'-0.56%"
$ws.Range("D16").Value = "'3.664"
$ws.Range("E16").Value = "'-1.66%"
$ws.Range("D17").Value = "'4.268"
$ws.Range("E17").Value = "'1.46%"
$ws.Range("D19").Value = "'0.3465"
$ws.Range("E19").Value = "'0.14%"
$ws.Range("D20").Value = "'0.1294"
$ws.Range("E20").Value = "'-3.01%"
$ws.Range("D21").Value = "'4.882"
$ws.Range("E21").Value = "'0.79%"
$ws.Range("E22").Value = "'1.49%"
$ws.Range("D23").Value = "'0.04320"
$ws.Range("E23").Value = "'-0.76%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-0.84%"
$ws.Range("D25").Value = "'0.004790"
$ws.Range("E25").Value = "'16.55%"
$ws.Range("D26").Value = "'0.0001296"
$ws.Range("E26").Value = "'-0.65%"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("E27").Value = "'-10.21%"
$ws.Range("D39").Value = "'0.02237"
$ws.Range("E39").Value = "'8.25%"
$ws.Range("D40").Value = "'0.05304"
$ws.Range("E40").Value = "'5.85%"
$ws.Range("D41").Value = "'0.007591"
$ws.Range("E41").Value = "'1.35%"
$ws.Range("D42").Value = "'0.009853"
$ws.Range("E42").Value = "'-2.74%"
$ws.Range("D43").Value = "'0.1380"
$ws.Range("E43").Value = "'2.69%"
$ws.Range("E44").Value = "'-1.64%"
$ws.Range("D45").Value = "'0.009786"
$ws.Range("E45").Value = "'9.83%"
$ws.Range("D46").Value = "'0.00006359"
$ws.Range("E46").Value = "'2.50%"
$ws.Range("E47").Value = "'-0.61%"
$ws.Range("D48").Value = "'0.001198"
$ws.Range("E48").Value = "'-25.29%"
$ws.Range("D49").Value = "'0.002759"
$ws.Range("E49").Value = "'-7.76%"
$ws.Range("D50").Value = "'0.00002094"
$ws.Range("E50").Value = "'-0.61%"
$ws.Range("E51").Value = "'-0.61%"
